$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 55
$ws1.Range("F4").Value = 1108
$ws1.Range("F7").Value = 591
$ws1.Range("F10").Value = 1420
$ws1.Range("F12").Value = 571
$ws1.Range("F13").Value = 1727
$ws1.Range("F14").Value = 1785
$ws1.Range("F17").Value = 1450
$ws1.Range("F18").Value = 279
$ws1.Range("F20").Value = 1181
$ws1.Range("F22").Value = 432
$ws1.Range("F23").Value = 61
$ws1.Range("F24").Value = 4666
$ws1.Range("F25").Value = 735
$ws1.Range("F26").Value = 566
$ws1.Range("F27").Value = 1615
$ws1.Range("F28").Value = 44
$ws1.Range("F29").Value = 90

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 45
$ws2.Range("F5").Value = 23
$ws2.Range("F9").Value = 51
$ws2.Range("F13").Value = 97
$ws2.Range("G3").Value = "不可售"

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 31

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 55
$ws4.Range("F4").Value = 31
$ws4.Range("F7").Value = 45
$ws4.Range("F8").Value = 23
$ws4.Range("F13").Value = 51
$ws4.Range("F15").Value = 1108
$ws4.Range("F18").Value = 591
$ws4.Range("F21").Value = 1420
$ws4.Range("F23").Value = 571
$ws4.Range("F24").Value = 1727
$ws4.Range("F25").Value = 1785
$ws4.Range("F28").Value = 1450
$ws4.Range("F29").Value = 279
$ws4.Range("F33").Value = 1181
$ws4.Range("F35").Value = 432
$ws4.Range("F36").Value = 61
$ws4.Range("F37").Value = 4666
$ws4.Range("F38").Value = 735
$ws4.Range("F39").Value = 566
$ws4.Range("F40").Value = 1615
$ws4.Range("F41").Value = 97
$ws4.Range("F43").Value = 44
$ws4.Range("F44").Value = 90
$ws4.Range("G6").Value = "不可售"
